$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in header: "Pedido Origen" -> "Pedido Origem"
$ws.Range("A1").Value = "Pedido Origem"

# Update "Data de Remessa" dates from 2025-08-30 to 2025-08-31 (serial 45899 -> 45900)
$ws.Range("O2").Value = 45900
$ws.Range("O3").Value = 45900

# Move the active selection to A2
$ws.Range("A2").Select() | Out-Null
